$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries recorded while adding "white" to the coding / analyses sheet.
$ws.Range("A95").Value = "WHITE"
$ws.Range("B95").Value = "yes"
$ws.Range("C95").Value = "yes                                   "

$ws.Range("A96").Value = "T0보다 CHECK쪽이 맞음"
$ws.Range("B96").Value = "yes"
$ws.Range("C96").Value = "yes"

# Leave the selection on the newly added rows, as in the source workbook.
$ws.Range("A95:C96").Select()
